# Fixes to the title slide:
#  1) Title shape: "... mer , der no php" -> "... mer , der"   (drop the "no php" part)
#  2) Subtitle shape: "Bruno Cezarcio" -> "Bruno Cezario"        (typo fix in the author's name)

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    for ($vi = 1; $vi -le $slide.Shapes.Count; $vi++) {
        $shape = $slide.Shapes.Item($vi)

        if (-not $shape.HasTextFrame) { continue }

        $tr = $shape.TextFrame.TextRange
        $text = $tr.Text

        $needle1 = " , der no php"
        $idx1 = $text.IndexOf($needle1)
        if ($idx1 -ge 0) {
            $sel = $tr.Characters($idx1 + 1, $needle1.Length)
            $sel.Text = " , der"
            $text = $tr.Text
        }

        $needle2 = "Bruno Cezarcio"
        $idx2 = $text.IndexOf($needle2)
        if ($idx2 -ge 0) {
            $sel = $tr.Characters($idx2 + 1, $needle2.Length)
            $sel.Text = "Bruno Cezario"
            $text = $tr.Text
        }
    }
}
